# Apply the edits from the commit:
#  - Sheet "6": rename "Tomi" (A14) to "P. Tomi", and add a new name
#    "H. Tomi" at the end of the list (A24). This sheet becomes the
#    active sheet/tab.
#  - Sheet "8": add a new name "Dávid" at the end of the list (A12).

$wb = $excel.ActiveWorkbook

# --- Sheet "6" ---
$ws6 = $wb.Worksheets.Item("6")
$ws6.Range("A14").Value = "P. Tomi"
$ws6.Range("A24").Value = "H. Tomi"
$ws6.Range("A24").Select() | Out-Null

# --- Sheet "8" ---
$ws8 = $wb.Worksheets.Item("8")
$ws8.Range("A12").Value = "Dávid"
$ws8.Range("A12").Select() | Out-Null

# Make sheet "6" the active sheet/tab (activeTab goes from 3 to 1)
$ws6.Activate() | Out-Null
